# cryptos.xlsx refresh: row-by-row Price / Volume(1h) updates, plus a
# name/link swap + refresh for the Quant <-> Aptos pair (rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" strings ("1.003", "0.7346", ...) parse as valid numbers,
# but the source data is inline text (dot-grouped, e.g. "29.792.24" for
# Bitcoin). Forcing text format before the write, then reverting the
# cell style to Normal afterwards, keeps the value textual without
# leaving a visible/applied number-format behind.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2  (Bitcoin)
$ws.Range('D2').Value = '29.792.24'
$ws.Range('E2').Value = '  -0.50%  '

# Row 3  (Ethereum)
$ws.Range('D3').Value = '1.863.09'
$ws.Range('E3').Value = '  -1.60%  '

# Row 4  (TetherUSD)
Set-TextValue 'D4' '1.003'
$ws.Range('E4').Value = '  +0.32%  '

# Row 5  (XRP)
Set-TextValue 'D5' '0.7346'
$ws.Range('E5').Value = '  -5.15%  '

# Row 6  (BNB)
Set-TextValue 'D6' '241.13'
$ws.Range('E6').Value = '  -1.14%  '

# Row 7  (USDC)
Set-TextValue 'D7' '1.003'
$ws.Range('E7').Value = '  +0.28%  '

# Row 8  (Cardano)
Set-TextValue 'D8' '0.3089'
$ws.Range('E8').Value = '  -1.50%  '

# Row 9  (Solana)
Set-TextValue 'D9' '24.52'
$ws.Range('E9').Value = '  -4.51%  '

# Row 10  (Dogecoin)
Set-TextValue 'D10' '0.07043'
$ws.Range('E10').Value = '  -4.10%  '

# Row 11  (TRON)
Set-TextValue 'D11' '0.08406'
$ws.Range('E11').Value = '  +4.28%  '

# Row 12  (Polygon)
Set-TextValue 'D12' '0.7465'
$ws.Range('E12').Value = '  -3.46%  '

# Row 13  (WrappedEther)
$ws.Range('D13').Value = '1.866.30'
$ws.Range('E13').Value = '  +1.25%  '

# Row 14  (Polkadot)
Set-TextValue 'D14' '5.313'
$ws.Range('E14').Value = '  -3.45%  '

# Row 15  (Litecoin)
Set-TextValue 'D15' '92.04'
$ws.Range('E15').Value = '  -2.18%  '

# Row 16  (WrappedBTC)
$ws.Range('D16').Value = '29.786.80'
$ws.Range('E16').Value = '  -0.40%  '

# Row 17  (Uniswap)
Set-TextValue 'D17' '6.036'
$ws.Range('E17').Value = '  -3.06%  '

# Row 18  (Avalanche)
Set-TextValue 'D18' '13.51'
$ws.Range('E18').Value = '  -3.61%  '

# Row 19  (BitcoinCash)
Set-TextValue 'D19' '239.57'
$ws.Range('E19').Value = '  -2.81%  '

# Row 20  (ShibaInu)
Set-TextValue 'D20' '0.000007772'
$ws.Range('E20').Value = '  -1.17%  '

# Row 21  (Dai)
Set-TextValue 'D21' '1.003'
$ws.Range('E21').Value = '  +0.26%  '

# Row 22  (WrappedliquidstakedEther2.0)
$ws.Range('D22').Value = '2.136.49'
$ws.Range('E22').Value = '  +1.21%  '

# Row 23  (BinanceUSD)
$ws.Range('E23').Value = '  +0.34%  '

# Row 24  (Chainlink)
Set-TextValue 'D24' '7.885'
$ws.Range('E24').Value = '  -3.31%  '

# Row 25  (Stellar)
Set-TextValue 'D25' '0.1560'
$ws.Range('E25').Value = '  -0.76%  '

# Row 26  (Cosmos)
Set-TextValue 'D26' '9.249'
$ws.Range('E26').Value = '  -2.09%  '

# Row 27  (Monero)
Set-TextValue 'D27' '162.02'
$ws.Range('E27').Value = '  -0.23%  '

# Row 28  (EthereumClassic)
Set-TextValue 'D28' '18.50'
$ws.Range('E28').Value = '  -1.42%  '

# Row 29  (LidoDAOToken)
Set-TextValue 'D29' '1.994'
$ws.Range('E29').Value = '  -1.58%  '

# Row 30  (Toncoin)
Set-TextValue 'D30' '1.488'
$ws.Range('E30').Value = '  +4.47%  '

# Row 31  (PancakeSwap)
Set-TextValue 'D31' '1.526'
$ws.Range('E31').Value = '  -0.97%  '

# Row 32  (Filecoin)
Set-TextValue 'D32' '4.435'
$ws.Range('E32').Value = '  -0.91%  '

# Row 33  (InternetComputer(DFINITY))
Set-TextValue 'D33' '4.124'
$ws.Range('E33').Value = '  +1.39%  '

# Row 34  (Hedera)
Set-TextValue 'D34' '0.05351'
$ws.Range('E34').Value = '  -3.88%  '

# Row 35  (ARBITRUM)
Set-TextValue 'D35' '1.225'
$ws.Range('E35').Value = '  -0.97%  '

# Row 36  (ImmutableX)
Set-TextValue 'D36' '0.7410'
$ws.Range('E36').Value = '  -1.18%  '

# Row 37  (Frax)
Set-TextValue 'D37' '1.003'
$ws.Range('E37').Value = '  +0.34%  '

# Row 38  (HuobiToken)
Set-TextValue 'D38' '2.697'
$ws.Range('E38').Value = '  +0.50%  '

# Row 39  (VeChain)
Set-TextValue 'D39' '0.01924'
$ws.Range('E39').Value = '  -0.33%  '

# Row 40  (MXToken)
Set-TextValue 'D40' '2.745'
$ws.Range('E40').Value = '  -1.57%  '

# Row 41  (TheSandbox)
Set-TextValue 'D41' '0.4410'
$ws.Range('E41').Value = '  -1.46%  '

# Row 42  (Maker)
$ws.Range('D42').Value = '1.096.33'
$ws.Range('E42').Value = '  -0.34%  '

# Row 43  (FraxShare)
Set-TextValue 'D43' '5.987'
$ws.Range('E43').Value = '  -0.28%  '

# Row 44  (Aave)
Set-TextValue 'D44' '71.58'
$ws.Range('E44').Value = '  -3.80%  '

# Row 45  (TrustWalletToken)
$ws.Range('E45').Value = '  +1.05%  '

# Row 46  (PaxDollar)
Set-TextValue 'D46' '1.003'
$ws.Range('E46').Value = '  +0.28%  '

# Row 47  (Quant)
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D47' '7.701'
$ws.Range('E47').Value = '  +2.09%  '

# Row 48  (Aptos)
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D48' '101.78'
$ws.Range('E48').Value = '  -0.70%  '

# Row 49  (RenderToken)
Set-TextValue 'D49' '1.823'
$ws.Range('E49').Value = '  -3.44%  '

# Row 50  (SynthetixNetwork)
Set-TextValue 'D50' '2.993'
$ws.Range('E50').Value = '  +0.06%  '

# Row 51  (RocketPoolETH)
$ws.Range('D51').Value = '2.034.31'
$ws.Range('E51').Value = '  +0.80%  '
